# Fix report acdc project & menu hrd
# - Update report period from Jan/23 to Dec/22 on both sheets
# - Remove the (now deleted) transaction-maker detail rows on sheet 1
# - Reset the total nominal to Rp0 and drop the stale disti rows on sheet 2

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Transaction Maker DCL" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Value = "Export Laporan Transaction Maker DCL Periode 01/Dec/22 - 31/Dec/22"
$ws1.Rows("3:4").Delete()

# --- Sheet 2: "Worksheet" ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Value = "Export Laporan Data DCL Periode 01/Dec/22 - 31/Dec/22"
$ws2.Range("B2").Value = "Rp0"
$ws2.Rows("5:6").Delete()
